$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
    $cell.ClearFormats()
}

# Row 2
Set-TextCell $ws "D2" "62.453.34"
Set-TextCell $ws "E2" "  +1.95%  "

# Row 3
Set-TextCell $ws "D3" "3.443.48"
Set-TextCell $ws "E3" "  +2.06%  "

# Row 4
Set-TextCell $ws "E4" "  +0.19%  "

# Row 5
Set-TextCell $ws "D5" "406.06"
Set-TextCell $ws "E5" "  -2.77%  "

# Row 6
Set-TextCell $ws "D6" "129.67"
Set-TextCell $ws "E6" "  +12.72%  "

# Row 7
Set-TextCell $ws "D7" "0.615"
Set-TextCell $ws "E7" "  +3.59%  "

# Row 8
Set-TextCell $ws "D8" "3.436.55"
Set-TextCell $ws "E8" "  +2.12%  "

# Row 9
Set-TextCell $ws "E9" "  +0.07%  "

# Row 10
Set-TextCell $ws "D10" "0.682"
Set-TextCell $ws "E10" "  +5.58%  "

# Row 11
Set-TextCell $ws "D11" "0.130"
Set-TextCell $ws "E11" "  +23.55%  "

# Row 12
Set-TextCell $ws "D12" "42.58"
Set-TextCell $ws "E12" "  +5.18%  "

# Row 13
Set-TextCell $ws "E13" "  -0.95%  "

# Row 14
Set-TextCell $ws "D14" "4.000.99"
Set-TextCell $ws "E14" "  +2.63%  "

# Row 15
Set-TextCell $ws "D15" "8.66"
Set-TextCell $ws "E15" "  +1.73%  "

# Row 16
Set-TextCell $ws "D16" "19.90"
Set-TextCell $ws "E16" "  -0.59%  "

# Row 17
Set-TextCell $ws "D17" "3.454.20"
Set-TextCell $ws "E17" "  +2.73%  "

# Row 18
Set-TextCell $ws "D18" "62.452.24"
Set-TextCell $ws "E18" "  +2.35%  "

# Row 19
Set-TextCell $ws "D19" "1.03"
Set-TextCell $ws "E19" "  -1.95%  "

# Row 20
Set-TextCell $ws "D20" "11.05"
Set-TextCell $ws "E20" "  +1.87%  "

# Row 21
Set-TextCell $ws "D21" "0.0000138"
Set-TextCell $ws "E21" "  +21.51%  "

# Row 22
Set-TextCell $ws "D22" "3.30"
Set-TextCell $ws "E22" "  -3.27%  "

# Row 23
Set-TextCell $ws "D23" "82.54"
Set-TextCell $ws "E23" "  +9.12%  "

# Row 24
Set-TextCell $ws "D24" "13.03"
Set-TextCell $ws "E24" "  -1.30%  "

# Row 25
Set-TextCell $ws "D25" "309.20"
Set-TextCell $ws "E25" "  +1.62%  "

# Row 26
Set-TextCell $ws "D26" "3.15"
Set-TextCell $ws "E26" "  -2.27%  "

# Row 27
Set-TextCell $ws "D27" "29.99"
Set-TextCell $ws "E27" "  +3.71%  "

# Row 28
Set-TextCell $ws "D28" "8.20"
Set-TextCell $ws "E28" "  +2.17%  "

# Row 29
Set-TextCell $ws "B29" "LEO"
Set-TextCell $ws "C29" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell $ws "D29" "4.37"
Set-TextCell $ws "E29" "  -2.61%  "

# Row 30
Set-TextCell $ws "B30" "RenderToken"
Set-TextCell $ws "C30" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws "D30" "7.54"
Set-TextCell $ws "E30" "  -2.51%  "

# Row 31
Set-TextCell $ws "B31" "Kaspa"
Set-TextCell $ws "C31" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws "D31" "0.176"
Set-TextCell $ws "E31" "  -2.08%  "

# Row 32
Set-TextCell $ws "D32" "0.116"
Set-TextCell $ws "E32" "  +0.68%  "

# Row 33
Set-TextCell $ws "D33" "43.26"
Set-TextCell $ws "E33" "  +7.39%  "

# Row 34
Set-TextCell $ws "D34" "11.68"
Set-TextCell $ws "E34" "  +0.83%  "

# Row 35
Set-TextCell $ws "D35" "2.57"
Set-TextCell $ws "E35" "  -1.99%  "

# Row 36
Set-TextCell $ws "D36" "0.998"
Set-TextCell $ws "E36" "  -0.11%  "

# Row 37
Set-TextCell $ws "D37" "0.0488"
Set-TextCell $ws "E37" "  -4.80%  "

# Row 38
Set-TextCell $ws "D38" "52.74"
Set-TextCell $ws "E38" "  +0.31%  "

# Row 39
Set-TextCell $ws "B39" "FirstDigitalUSD"
Set-TextCell $ws "C39" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell $ws "D39" "1.00"
Set-TextCell $ws "E39" "  +0.25%  "

# Row 40
Set-TextCell $ws "B40" "LidoDAOToken"
Set-TextCell $ws "C40" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell $ws "D40" "3.48"
Set-TextCell $ws "E40" "  +1.25%  "

# Row 41
Set-TextCell $ws "D41" "3.00"
Set-TextCell $ws "E41" "  -4.02%  "

# Row 42
Set-TextCell $ws "D42" "0.127"
Set-TextCell $ws "E42" "  +2.70%  "

# Row 43
Set-TextCell $ws "B43" "Monero"
Set-TextCell $ws "C43" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws "D43" "137.54"
Set-TextCell $ws "E43" "  +0.05%  "

# Row 44
Set-TextCell $ws "B44" "ARBITRUM"
Set-TextCell $ws "C44" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws "D44" "1.98"
Set-TextCell $ws "E44" "  +2.06%  "

# Row 45
Set-TextCell $ws "D45" "17.19"
Set-TextCell $ws "E45" "  +0.28%  "

# Row 46
Set-TextCell $ws "D46" "0.285"
Set-TextCell $ws "E46" "  -2.83%  "

# Row 47
Set-TextCell $ws "D47" "3.92"
Set-TextCell $ws "E47" "  -2.24%  "

# Row 48
Set-TextCell $ws "E48" "  -0.14%  "

# Row 49
Set-TextCell $ws "D49" "22.04"
Set-TextCell $ws "E49" "  -2.94%  "

# Row 50
Set-TextCell $ws "D50" "3.791.44"
Set-TextCell $ws "E50" "  +2.68%  "

# Row 51
Set-TextCell $ws "D51" "2.164.61"
Set-TextCell $ws "E51" "  -0.53%  "
